# Automatische test-sync: 2025-06-26 21:16:50
# Adds a new Logs row (row 17) for an incoming mail, extends the
# conditionalFormatting ranges to cover it, and bumps the Dashboard
# "Bestelling / Levering" count.

$wb = $excel.ActiveWorkbook

# --- Logs sheet: append the new mail row -----------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A17").Value = "Kun je 5 liter koelvloeistof bestellen?"
$logs.Range("B17").Value = "MailMind Test <mailmind.test@zohomail.eu>"
$logs.Range("C17").Value = "Hoi Johan,`n Zou je 5 liter koelvloeistof kunnen bestellen voor de werkplaats?`n Laat me even weten of dit lukt.`nGroet,`n Rick`nSent using {0}"
$logs.Range("D17").Value = "Bestelling / Levering"
$logs.Range("E17").Value = "Bedankt voor je bericht. Ik neem dit z.s.m. in behandeling."
$logs.Range("F17").Value = "2025-06-26 21:16:24"
$logs.Range("G17").Value = "Ja"
$logs.Range("H17").Value = "Ja"
$logs.Range("I17").Value = "Nee"

# The multi-line content in C17 would otherwise make the engine stamp an
# explicit row height (ht/customHeight); re-auto-fit so row 17 serializes
# the same way the other (height-less) rows do.
$logs.Rows.Item(17).AutoFit()

# Extend the conditional-formatting sqrefs (D/G/H/I, rows 2-16) so they
# also cover the newly added row 17.
$dFc = $logs.Range("D2:D16").FormatConditions
for ($i = 1; $i -le $dFc.Count; $i++) {
    $dFc.Item($i).ModifyAppliesToRange($logs.Range("D2:D17"))
}

$gFc = $logs.Range("G2:G16").FormatConditions
for ($i = 1; $i -le $gFc.Count; $i++) {
    $gFc.Item($i).ModifyAppliesToRange($logs.Range("G2:G17"))
}

$hFc = $logs.Range("H2:H16").FormatConditions
for ($i = 1; $i -le $hFc.Count; $i++) {
    $hFc.Item($i).ModifyAppliesToRange($logs.Range("H2:H17"))
}

$iFc = $logs.Range("I2:I16").FormatConditions
for ($i = 1; $i -le $iFc.Count; $i++) {
    $iFc.Item($i).ModifyAppliesToRange($logs.Range("I2:I17"))
}

# --- Dashboard sheet: bump the "Bestelling / Levering" count ---------
$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Range("B2").Value = 12
